$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2243
$ws.Range("J18").Value = 2759.8
$ws.Range("L18").Value = 2759.8
$ws.Range("N18").Value = -3327.8

$ws.Range("H33").Value = 335.7647
$ws.Range("I33").Value = 324.69232
$ws.Range("K33").Value = 324.69232
$ws.Range("M33").Value = -95.69232

$ws.Range("H129").Value = 4271.778
$ws.Range("J129").Value = 4570
$ws.Range("L129").Value = 13710
$ws.Range("N129").Value = -23710

$ws.Range("H132").Value = 5343.2173
$ws.Range("I132").Value = 5376.095
$ws.Range("K132").Value = 16128.285
$ws.Range("M132").Value = -13598.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1765
$ws.Range("I21").Value = 1597.5
$ws.Range("J21").Value = 2100
$ws.Range("K21").Value = 1597.5
$ws.Range("L21").Value = 2100
$ws.Range("M21").Value = -1223.5
$ws.Range("N21").Value = -2848

$ws.Range("H63").Value = 3750.5527
$ws.Range("I63").Value = 2394.375
$ws.Range("K63").Value = 2394.375
$ws.Range("M63").Value = -1708.375

$ws.Range("H66").Value = 3750.5527
$ws.Range("I66").Value = 2394.375
$ws.Range("K66").Value = 11971.875
$ws.Range("M66").Value = -8539.875

$ws.Range("H74").Value = 1922.25
$ws.Range("I74").Value = 1654.3334
$ws.Range("K74").Value = 1654.3334
$ws.Range("M74").Value = -780.3334

$ws.Range("H77").Value = 1922.25
$ws.Range("I77").Value = 1654.3334
$ws.Range("K77").Value = 8271.666999999999
$ws.Range("M77").Value = -3903.666999999999

$ws.Range("H80").Value = 101994
$ws.Range("J80").Value = 101994
$ws.Range("L80").Value = 101994
$ws.Range("N80").Value = -103990

$ws.Range("H83").Value = 101994
$ws.Range("J83").Value = 101994
$ws.Range("L83").Value = 305982
$ws.Range("N83").Value = -315966

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H135").Value = 49423
$ws.Range("J135").Value = 49423
$ws.Range("L135").Value = 49423
$ws.Range("N135").Value = -59563

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 16999.2
$ws.Range("J103").Value = 16999.2
$ws.Range("L103").Value = 16999.2
$ws.Range("N103").Value = -19343.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7277.65
$ws.Range("I22").Value = 1823.5
$ws.Range("J22").Value = 10913.75
$ws.Range("K22").Value = 1823.5
$ws.Range("L22").Value = 10913.75
$ws.Range("M22").Value = -1473.5
$ws.Range("N22").Value = -11613.75

$ws.Range("H31").Value = 3720
$ws.Range("I31").Value = 3354.4285
$ws.Range("K31").Value = 3354.4285
$ws.Range("M31").Value = -3059.4285

$ws.Range("H34").Value = 3720
$ws.Range("I34").Value = 3354.4285
$ws.Range("K34").Value = 3354.4285
$ws.Range("M34").Value = -3152.4285

$ws.Range("H134").Value = 3964.2
$ws.Range("I134").Value = 3849.111
$ws.Range("K134").Value = 11547.333
$ws.Range("M134").Value = -9012.332999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 286.2857
$ws.Range("J97").Value = 210.25
$ws.Range("L97").Value = 630.75
$ws.Range("N97").Value = -1622.75

$ws.Range("H114").Value = 1807
$ws.Range("I114").Value = 949.4
$ws.Range("J114").Value = 2419.5715
$ws.Range("K114").Value = 2848.2
$ws.Range("L114").Value = 7258.7145
$ws.Range("M114").Value = 405.8000000000002
$ws.Range("N114").Value = -13766.7145

$ws.Range("H125").Value = 27501.5
$ws.Range("J125").Value = 34998.5
$ws.Range("L125").Value = 104995.5
$ws.Range("N125").Value = -114835.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1918.4286
$ws.Range("I43").Value = 1423
$ws.Range("J43").Value = 3157
$ws.Range("K43").Value = 1423
$ws.Range("L43").Value = 3157
$ws.Range("M43").Value = -1272
$ws.Range("N43").Value = -3459

$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 1000
$ws.Range("M80").Value = -2
$ws.Range("N80").Value = -2996

$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 5000
$ws.Range("L83").Value = 5000
$ws.Range("M83").Value = -8
$ws.Range("N83").Value = -14984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 8
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H43").Value = 248000
$ws.Range("I43").Value = 18272.727
$ws.Range("J43").Value = 669166.7
$ws.Range("K43").Value = 18272.727
$ws.Range("L43").Value = 669166.7
$ws.Range("M43").Value = -18079.727
$ws.Range("N43").Value = -669552.7

$ws.Range("H53").Value = 5000
$ws.Range("I53").Value = 5000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 5000
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -4482
$ws.Range("N53").ClearContents()

$ws.Range("H132").Value = 2323.9473
$ws.Range("I132").Value = 2250.4119
$ws.Range("K132").Value = 6751.2357
$ws.Range("M132").Value = -4221.2357

$ws.Range("H136").Value = 2422.7334
$ws.Range("I136").Value = 2549.3845
$ws.Range("J136").Value = 1599.5
$ws.Range("K136").Value = 7648.1535
$ws.Range("L136").Value = 4798.5
$ws.Range("M136").Value = -5098.1535
$ws.Range("N136").Value = -9898.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 10427.75
$ws.Range("I38").Value = 7555.5
$ws.Range("J38").Value = 13300
$ws.Range("K38").Value = 7555.5
$ws.Range("L38").Value = 13300
$ws.Range("M38").Value = -7082.5
$ws.Range("N38").Value = -14246

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H81").Value = 4025.9333
$ws.Range("J81").Value = 5958.6
$ws.Range("L81").Value = 11917.2
$ws.Range("N81").Value = -14039.2

$ws.Range("H84").Value = 4025.9333
$ws.Range("J84").Value = 5958.6
$ws.Range("L84").Value = 59586
$ws.Range("N84").Value = -70194

$ws.Range("H97").Value = 14714
$ws.Range("J97").Value = 14714
$ws.Range("L97").Value = 14714
$ws.Range("N97").Value = -16696
